$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells hold numeric-looking values stored as text (inline strings).
# Set NumberFormat to Text before assigning so Excel keeps the literal
# string (preserving trailing zeros / decimals), then clear the format
# back to the default so no stray cell style is introduced.
$cellValues = @{
    "C21" = "49"
    "E21" = "141171.00"
    "C24" = "519"
    "D24" = "430"
    "E24" = "4439957.92"
    "C28" = "58"
    "E28" = "350148.39"
    "C32" = "100"
    "E32" = "1219245.89"
    "C84" = "846"
    "E84" = "7693397.41"
    "C92" = "132"
    "E92" = "1221799.93"
    "C130" = "1162"
    "E130" = "9866904.14"
    "C134" = "170"
    "E134" = "1695798.92"
    "C147" = "5096"
    "E147" = "35690046.99"
    "C151" = "1642"
    "E151" = "6317634.03"
    "C155" = "861"
    "E155" = "5245048.93"
    "C163" = "16"
    "D163" = "16"
    "E163" = "170500.00"
}

foreach ($cellRef in $cellValues.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $cellValues[$cellRef]
    $range.ClearFormats()
}
